$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The shared status string "Ready for handoff" is used by the Overview
# rollup cells (B3/C3) as well as the per-language sheets (C3). Update all
# occurrences so the report reflects that the handback transform failed.
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file-name mismatch error detail for the
# failed file on each language sheet.
$wsZhCn.Range("K3").Value = "Handback file name: qcuwzh1e.iwh is different with handoff file name: 897f35b8-fa64-40e2-9cb9-29cb7ea1defe.0a7ad118ac28c9a292f3e200031a1d0ea68e0000.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: qcuwzh1e.iwh is different with handoff file name: 897f35b8-fa64-40e2-9cb9-29cb7ea1defe.0a7ad118ac28c9a292f3e200031a1d0ea68e0000.de-de."
